# Add a new "Non-Functional Requirements" slide right before the final
# (image) slide. The final slide currently sits at index 16 (id 267);
# inserting a new Title+Content slide at index 16 pushes it to index 17,
# matching the target sldIdLst ordering (…, 269, 273[new], 267).

$p = $ppt.ActivePresentation

# ppLayoutText (2) == "Title and Content" (slideLayout2.xml), the layout
# used by every other content slide in this deck.
$s = $p.Slides.Add(16, 2)

# --- Title placeholder -----------------------------------------------
$title = $s.Shapes.Item(1)

$title.Left   = 822960  / 12700.0
$title.Top    = 759656  / 12700.0
$title.Width  = 7520940 / 12700.0
$title.Height = 548640  / 12700.0

$titleTr = $title.TextFrame.TextRange
$titleTr.Text = "Non-Functional Requirements"
$titleTr.InsertAfter([char]0x0B) | Out-Null

# Centre only the first (visible) line, matching the authored slide.
$titleTr.Paragraphs(1, 1).ParagraphFormat.Alignment = 2

# --- Content placeholder ----------------------------------------------
$body = $s.Shapes.Item(2)

$body.Left   = 822960  / 12700.0
$body.Top    = 1677403 / 12700.0
$body.Width  = 7520940 / 12700.0
$body.Height = 3579849 / 12700.0

$dot = [char]0xB7

$bullets = @(
    "$dot        8 letter username for users",
    "$dot        At least 5 users able to be logged in to the system",
    "$dot        Keep track of certain amount of orders",
    "$dot        Supports popular browsers (IE, Chrome, Firefox)",
    "$dot        Passwords (At least 1 Capital, 1 number, 8 characters)",
    "$dot        Recovery options for Data",
    "$dot        Maintenance",
    "$dot        Hidden Items, separate database",
    "$dot        Account Security (Usernames and Passwords)"
)

# Build the bullets one paragraph at a time via InsertAfter (rather than
# one "Text = a`ra`rb`r..." assignment) so every run gets its own
# per-language run properties, matching authored PowerPoint output.
$bodyTr = $body.TextFrame.TextRange
$bodyTr.Text = $bullets[0]
for ($i = 1; $i -lt $bullets.Length; $i++) {
    $bodyTr.InsertAfter([char]0x0D + $bullets[$i]) | Out-Null
}
